$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-08 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-09 Monday", 2) | Out-Null
$d.Content.Find.Execute("25×47=1175", $true, $false, $false, $false, $false, $true, 1, $false, "50×90=4500", 2) | Out-Null
$d.Content.Find.Execute("41×47=1927", $true, $false, $false, $false, $false, $true, 1, $false, "47×18=846", 2) | Out-Null
$d.Content.Find.Execute("23×38=874", $true, $false, $false, $false, $false, $true, 1, $false, "84×34=2856", 2) | Out-Null
$d.Content.Find.Execute("86×96=8256", $true, $false, $false, $false, $false, $true, 1, $false, "76×80=6080", 2) | Out-Null
$d.Content.Find.Execute("71×20=1420", $true, $false, $false, $false, $false, $true, 1, $false, "84×68=5712", 2) | Out-Null
$d.Content.Find.Execute("65×88=5720", $true, $false, $false, $false, $false, $true, 1, $false, "60×39=2340", 2) | Out-Null
$d.Content.Find.Execute("89×82=7298", $true, $false, $false, $false, $false, $true, 1, $false, "39×76=2964", 2) | Out-Null
$d.Content.Find.Execute("79×74=5846", $true, $false, $false, $false, $false, $true, 1, $false, "60×87=5220", 2) | Out-Null
$d.Content.Find.Execute("59×66=3894", $true, $false, $false, $false, $false, $true, 1, $false, "53×79=4187", 2) | Out-Null
$d.Content.Find.Execute("31×49=1519", $true, $false, $false, $false, $false, $true, 1, $false, "24×58=1392", 2) | Out-Null
$d.Content.Find.Execute("44×46=2024", $true, $false, $false, $false, $false, $true, 1, $false, "19×50=950", 2) | Out-Null
$d.Content.Find.Execute("36×58=2088", $true, $false, $false, $false, $false, $true, 1, $false, "80×64=5120", 2) | Out-Null
$d.Content.Find.Execute("87×93=8091", $true, $false, $false, $false, $false, $true, 1, $false, "17×13=221", 2) | Out-Null
$d.Content.Find.Execute("56×68=3808", $true, $false, $false, $false, $false, $true, 1, $false, "51×94=4794", 2) | Out-Null
$d.Content.Find.Execute("68×12=816", $true, $false, $false, $false, $false, $true, 1, $false, "74×47=3478", 2) | Out-Null
$d.Content.Find.Execute("31×77=2387", $true, $false, $false, $false, $false, $true, 1, $false, "12×65=780", 2) | Out-Null
$d.Content.Find.Execute("31×63=1953", $true, $false, $false, $false, $false, $true, 1, $false, "57×74=4218", 2) | Out-Null
$d.Content.Find.Execute("90×35=3150", $true, $false, $false, $false, $false, $true, 1, $false, "51×71=3621", 2) | Out-Null
$d.Content.Find.Execute("28×27=756", $true, $false, $false, $false, $false, $true, 1, $false, "87×46=4002", 2) | Out-Null
$d.Content.Find.Execute("20×74=1480", $true, $false, $false, $false, $false, $true, 1, $false, "21×98=2058", 2) | Out-Null
$d.Content.Find.Execute("66×88=5808", $true, $false, $false, $false, $false, $true, 1, $false, "82×26=2132", 2) | Out-Null
$d.Content.Find.Execute("52×52=2704", $true, $false, $false, $false, $false, $true, 1, $false, "53×74=3922", 2) | Out-Null
$d.Content.Find.Execute("44×17=748", $true, $false, $false, $false, $false, $true, 1, $false, "62×41=2542", 2) | Out-Null
$d.Content.Find.Execute("88×32=2816", $true, $false, $false, $false, $false, $true, 1, $false, "76×45=3420", 2) | Out-Null
$d.Content.Find.Execute("54×80=4320", $true, $false, $false, $false, $false, $true, 1, $false, "94×23=2162", 2) | Out-Null
